{"js": "// The document had two small textual edits:\n//  1. \"Semester: Fall 202\" -> \"Semester: Fall 2020\" (append the digit \"0\")\n//  2. \"...more straight forward to conduct...\" ->\n//     \"...more straight forward ways to conduct...\" (insert the word \"ways \")\n//\n// (The rest of the canonical-XML diff is purely incidental run-splitting /\n// run-merging noise introduced by Word's own re-serialization \u2014 the visible\n// paragraph text is identical before and after except for these two\n// insertions, matching the commit message \"adding a word to the file\".)\n\nconst body = context.document.body;\n\n// 1) \"Semester: Fall 202\" -> \"Semester: Fall 2020\"\nconst semesterResults = body.search(\"Semester: Fall 202\", { matchCase: true });\nsemesterResults.load(\"text\");\nawait context.sync();\n\nif (semesterResults.items.length > 0) {\n  semesterResults.items[0].insertText(\"0\", Word.InsertLocation.end);\n  await context.sync();\n}\n\n// 2) Insert the word \"ways \" before \"to conduct classification analysis\"\nconst waysResults = body.search(\"to conduct classification analysis\", { matchCase: true });\nwaysResults.load(\"text\");\nawait context.sync();\n\nif (waysResults.items.length > 0) {\n  waysResults.items[0].insertText(\"ways \", Word.InsertLocation.before);\n  await context.sync();\n}\n", "ps1": "# The document had two small textual edits:\n#  1. \"Semester: Fall 202\" -> \"Semester: Fall 2020\" (append the digit \"0\")\n#  2. \"...more straight forward to conduct...\" ->\n#     \"...more straight forward ways to conduct...\" (insert the word \"ways \")\n#\n# (The rest of the canonical-XML diff is purely incidental run-splitting /\n# run-merging noise introduced by Word's own re-serialization -- the visible\n# paragraph text is identical before and after except for these two\n# insertions, matching the commit message \"adding a word to the file\".)\n\n$d = $word.ActiveDocument\n\n# 1) \"Semester: Fall 202\" -> \"Semester: Fall 2020\"\n$rng1 = $d.Content\n$rng1.Find.Execute(\"Semester: Fall 202\") | Out-Null\nif ($rng1.Find.Found) {\n    $rng1.Collapse(0)   # wdCollapseEnd\n    $rng1.InsertAfter(\"0\")\n}\n\n# 2) Insert the word \"ways \" before \"to conduct classification analysis\"\n$rng2 = $d.Content\n$rng2.Find.Execute(\"to conduct classification analysis\") | Out-Null\nif ($rng2.Find.Found) {\n    $rng2.Collapse(1)   # wdCollapseStart\n    $rng2.InsertBefore(\"ways \")\n}\n"}
